$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; E=2; G=11.66862; H=23.33724; I=0.1418293637844402; J=0.1026592641785845; K=2; M=163.8679735; N=327.735947; O=0.4353980439170849; P=0.3519187842322066; Q=1912.11311294157; R=7648.452451766279; S=0.0617522275617499; T=0.03612772343990039}
    @{Row=3; E=2; G=11.66862; H=23.33724; I=0.1418293637844402; J=0.1026592641785845; K=3; M=16.807086; N=50.42125799999999; O=0.04465651348490144; P=0.05414171981207304; Q=196.11549984132; R=1176.69299904792; S=0.006333604896394845; T=0.005558149117270509}
    @{Row=4; E=2; G=11.66862; H=23.33724; I=0.1418293637844402; J=0.1026592641785845; K=3; M=31.35339366666667; N=94.060181; O=0.08330612737228353; P=0.1010006526448602; Q=365.85083640674; R=2195.10501844044; S=0.01181525504455651; T=0.01036865268207815}
    @{Row=5; E=2; G=11.66862; H=23.33724; I=0.1418293637844402; J=0.1026592641785845; K=3; M=35.54054833333333; N=106.621645; O=0.09443141874256437; P=0.1144889964762942; Q=414.7091530932999; R=2488.2549185598; S=0.01339314804151997; T=0.01175335613480091}
    @{Row=6; E=2; G=11.66862; H=23.33724; I=0.1418293637844402; J=0.1026592641785845; K=3; M=94.85468300000001; N=284.564049; O=0.2520293779954212; P=0.3055613370365933; Q=1106.82325114746; R=6640.93950688476; S=0.03574516633607878; T=0.03136870202160113}
    @{Row=7; E=2; G=11.66862; H=23.33724; I=0.1418293637844402; J=0.1026592641785845; K=2; M=33.9399115; N=67.879823; O=0.0901785184877445; P=0.07288850979797275; Q=396.03193012713; R=1584.12772050852; S=0.01278996190414018; T=0.007482680782933431}
    @{Row=8; E=3; G=15.12099333333333; H=45.36298; I=0.1837921591632454; J=0.1995493103618014; K=2; M=163.8679735; N=327.735947; O=0.4353980439170849; P=0.3519187842322066; Q=2477.846534840343; R=14867.07920904206; S=0.08002274658697456; T=0.07022515069690043}
    @{Row=9; E=3; G=15.12099333333333; H=45.36298; I=0.1837921591632454; J=0.1995493103618014; K=3; M=16.807086; N=50.42125799999999; O=0.04465651348490144; P=0.05414171981207304; Q=254.13983535876; R=2287.25851822884; S=0.008207517034092619; T=0.01080394285030106}
    @{Row=10; E=3; G=15.12099333333333; H=45.36298; I=0.1837921591632454; J=0.1995493103618014; K=3; M=31.35339366666667; N=94.060181; O=0.08330612737228353; P=0.1010006526448602; Q=474.0944566110422; R=4266.85010949938; S=0.01531101302128033; T=0.0201546105813737}
    @{Row=11; E=3; G=15.12099333333333; H=45.36298; I=0.1837921591632454; J=0.1995493103618014; K=3; M=35.54054833333333; N=106.621645; O=0.09443141874256437; P=0.1144889964762942; Q=537.4083944113444; R=4836.675549702099; S=0.01735575434354446; T=0.02284620029085921}
    @{Row=12; E=3; G=15.12099333333333; H=45.36298; I=0.1837921591632454; J=0.1995493103618014; K=3; M=94.85468300000001; N=284.564049; O=0.2520293779954212; P=0.3055613370365933; Q=1434.297029278447; R=12908.67326350602; S=0.0463210235543482; T=0.06097455407888216}
    @{Row=13; E=3; G=15.12099333333333; H=45.36298; I=0.1837921591632454; J=0.1995493103618014; K=2; M=33.9399115; N=67.879823; O=0.0901785184877445; P=0.07288850979797275; Q=513.2051755254233; R=3079.23105315254; S=0.0165741046230052; T=0.01454485186348487}
    @{Row=14; E=3; G=15.30100133333333; H=45.903004; I=0.1859801145612367; J=0.2019248469067731; K=2; M=163.8679735; N=327.735947; O=0.4353980439170849; P=0.3519187842322066; Q=2507.344081014131; R=15044.06448608479; S=0.08097537808743781; T=0.07106114662970606}
    @{Row=15; E=3; G=15.30100133333333; H=45.903004; I=0.1859801145612367; J=0.2019248469067731; K=3; M=16.807086; N=50.42125799999999; O=0.04465651348490144; P=0.05414171981207304; Q=257.165245295448; R=2314.487207659031; S=0.00830522349382738; T=0.01093255848432225}
    @{Row=16; E=3; G=15.30100133333333; H=45.903004; I=0.1859801145612367; J=0.2019248469067731; K=3; M=31.35339366666667; N=94.060181; O=0.08330612737228353; P=0.1010006526448602; Q=479.7383182981915; R=4317.644864683723; S=0.01549328311235027; T=0.02039454132279756}
    @{Row=17; E=3; G=15.30100133333333; H=45.903004; I=0.1859801145612367; J=0.2019248469067731; K=3; M=35.54054833333333; N=106.621645; O=0.09443141874256437; P=0.1144889964762942; Q=543.8059774357309; R=4894.253796921579; S=0.01756236607592224; T=0.02311817308598579}
    @{Row=18; E=3; G=15.30100133333333; H=45.903004; I=0.1859801145612367; J=0.2019248469067731; K=3; M=94.85468300000001; N=284.564049; O=0.2520293779954212; P=0.3055613370365933; Q=1451.371631055911; R=13062.3446795032; S=0.04687245259238567; T=0.06170042620174301}
    @{Row=19; E=3; G=15.30100133333333; H=45.903004; I=0.1859801145612367; J=0.2019248469067731; K=2; M=33.9399115; N=67.879823; O=0.0901785184877445; P=0.07288850979797275; Q=519.3146311147153; R=3115.887786688292; S=0.01677141119931332; T=0.01471800118221848}
    @{Row=20; E=3; G=24.037621; H=72.112863; I=0.2921716958236277; J=0.3172206076378816; K=2; M=163.8679735; N=327.735947; O=0.4353980439170849; P=0.3519187842322066; Q=3938.996241031044; R=23633.97744618626; S=0.127210984849545; T=0.1116358905733251}
    @{Row=21; E=3; G=24.037621; H=72.112863; I=0.2921716958236277; J=0.3172206076378816; K=3; M=16.807086; N=50.42125799999999; O=0.04465651348490144; P=0.05414171981207304; Q=404.002363382406; R=3636.021270441654; S=0.01304736927445435; T=0.01717486925734574}
    @{Row=22; E=3; G=24.037621; H=72.112863; I=0.2921716958236277; J=0.3172206076378816; K=3; M=31.35339366666667; N=94.060181; O=0.08330612737228353; P=0.1010006526448602; Q=753.6609940231336; R=6782.948946208204; S=0.02433969250685921; T=0.03203948840382515}
    @{Row=23; E=3; G=24.037621; H=72.112863; I=0.2921716958236277; J=0.3172206076378816; K=3; M=35.54054833333333; N=106.621645; O=0.09443141874256437; P=0.1144889964762942; Q=854.3102309688483; R=7688.792078719634; S=0.02759018775304614; T=0.03631826903006131}
    @{Row=24; E=3; G=24.037621; H=72.112863; I=0.2921716958236277; J=0.3172206076378816; K=3; M=94.85468300000001; N=284.564049; O=0.2520293779954212; P=0.3055613370365933; Q=2280.080920029143; R=20520.72828026229; S=0.07363585076629631; T=0.09693035300539164}
    @{Row=25; E=3; G=24.037621; H=72.112863; I=0.2921716958236277; J=0.3172206076378816; K=2; M=33.9399115; N=67.879823; O=0.0901785184877445; P=0.07288850979797275; Q=815.8347294105416; R=4895.008376463249; S=0.02634761067342668; T=0.0231217373679326}
    @{Row=26; E=3; G=8.323066666666668; H=24.9692; I=0.101164940675831; J=0.1098381684864154; K=2; M=163.8679735; N=327.735947; O=0.4353980439170849; P=0.3519187842322066; Q=1363.884067972067; R=8183.3044078324; S=0.04404701728324473; T=0.03865411471603158}
    @{Row=27; E=3; G=8.323066666666668; H=24.9692; I=0.101164940675831; J=0.1098381684864154; K=3; M=16.807086; N=50.42125799999999; O=0.04465651348490144; P=0.05414171981207304; Q=139.8864972504; R=1258.9784752536; S=0.0045176735374895; T=0.005946827342862774}
    @{Row=28; E=3; G=8.323066666666668; H=24.9692; I=0.101164940675831; J=0.1098381684864154; K=3; M=31.35339366666667; N=94.060181; O=0.08330612737228353; P=0.1010006526448602; Q=260.9563857139111; R=2348.6074714252; S=0.008427659433550281; T=0.01109372670244407}
    @{Row=29; E=3; G=8.323066666666668; H=24.9692; I=0.101164940675831; J=0.1098381684864154; K=3; M=35.54054833333333; N=106.621645; O=0.09443141874256437; P=0.1144889964762942; Q=295.8063531482222; R=2662.257178334; S=0.009553148875026077; T=0.01257526168480382}
    @{Row=30; E=3; G=8.323066666666668; H=24.9692; I=0.101164940675831; J=0.1098381684864154; K=3; M=94.85468300000001; N=284.564049; O=0.2520293779954212; P=0.3055613370365933; Q=789.4818502545335; R=7105.336652290801; S=0.02549653707347337; T=0.0335622976203597}
    @{Row=31; E=3; G=8.323066666666668; H=24.9692; I=0.101164940675831; J=0.1098381684864154; K=2; M=33.9399115; N=67.879823; O=0.0901785184877445; P=0.07288850979797275; Q=282.4841460752667; R=1694.9048764516; S=0.009122904473046999; T=0.008005940419913472}
    @{Row=32; E=2; G=7.8209415; H=15.641883; I=0.09506172599161901; J=0.06880780242854384; K=2; M=163.8679735; N=327.735947; O=0.4353980439170849; P=0.3519187842322066; Q=1281.60183446705; R=5126.407337868201; S=0.04138968954813282; T=0.02421475817634302}
    @{Row=33; E=2; G=7.8209415; H=15.641883; I=0.09506172599161901; J=0.06880780242854384; K=3; M=16.807086; N=50.42125799999999; O=0.04465651348490144; P=0.05414171981207304; Q=131.447236391469; R=788.6834183488139; S=0.00424512524864274; T=0.003725372759970699}
    @{Row=34; E=2; G=7.8209415; H=15.641883; I=0.09506172599161901; J=0.06880780242854384; K=3; M=31.35339366666667; N=94.060181; O=0.08330612737228353; P=0.1010006526448602; Q=245.2130576934705; R=1471.278346160823; S=0.00791922425368693; T=0.006949632952341521}
    @{Row=35; E=2; G=7.8209415; H=15.641883; I=0.09506172599161901; J=0.06880780242854384; K=3; M=35.54054833333333; N=106.621645; O=0.09443141874256437; P=0.1144889964762942; Q=277.9605493929225; R=1667.763296357535; S=0.008976813653505489; T=0.0078777362497831}
    @{Row=36; E=2; G=7.8209415; H=15.641883; I=0.09506172599161901; J=0.06880780242854384; K=3; M=94.85468300000001; N=284.564049; O=0.2520293779954212; P=0.3055613370365933; Q=741.8529267440446; R=4451.117560464268; S=0.02395834767283891; T=0.02102500410861561}
    @{Row=37; E=2; G=7.8209415; H=15.641883; I=0.09506172599161901; J=0.06880780242854384; K=2; M=33.9399115; N=67.879823; O=0.0901785184877445; P=0.07288850979797275; Q=265.4420623566772; R=1061.768249426709; S=0.008572525614812117; T=0.00501529818148989}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = $item.I
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 14).Value = $item.N
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
    $ws.Cells.Item($r, 17).Value = $item.Q
    $ws.Cells.Item($r, 18).Value = $item.R
    $ws.Cells.Item($r, 19).Value = $item.S
    $ws.Cells.Item($r, 20).Value = $item.T
}
